# Scheduled runner update: refresh market-price-derived columns (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the Goblin_Profits leve tables for each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1000423.3
$ws.Range("I2").Value = 1250372.6
$ws.Range("K2").Value = 1250372.6
$ws.Range("M2").Value = -1250259.6
$ws.Range("H15").Value = 453.34146
$ws.Range("I15").Value = 453.34146
$ws.Range("K15").Value = 1360.02438
$ws.Range("M15").Value = -1191.02438
$ws.Range("H33").Value = 1326.1
$ws.Range("I33").Value = 280
$ws.Range("K33").Value = 280
$ws.Range("M33").Value = -51
$ws.Range("H86").Value = 1811.3077
$ws.Range("I86").Value = 1215.2858
$ws.Range("J86").Value = 2506.6667
$ws.Range("K86").Value = 1215.2858
$ws.Range("L86").Value = 2506.6667
$ws.Range("M86").Value = -92.28580000000011
$ws.Range("N86").Value = -4752.6667
$ws.Range("H89").Value = 1811.3077
$ws.Range("I89").Value = 1215.2858
$ws.Range("J89").Value = 2506.6667
$ws.Range("K89").Value = 6076.429
$ws.Range("L89").Value = 12533.3335
$ws.Range("M89").Value = -460.4290000000001
$ws.Range("N89").Value = -23765.3335
$ws.Range("H97").Value = 5688
$ws.Range("J97").Value = 5688
$ws.Range("L97").Value = 17064
$ws.Range("N97").Value = -18056
$ws.Range("H103").Value = 1141.421
$ws.Range("J103").Value = 1123.52
$ws.Range("L103").Value = 3370.56
$ws.Range("N103").Value = -4542.559999999999
$ws.Range("H112").Value = 2023.0714
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 2165.7273
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 6497.1819
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -8713.1819
$ws.Range("H138").Value = 3459.307
$ws.Range("I138").Value = 2189.5417
$ws.Range("J138").Value = 3935.4688
$ws.Range("K138").Value = 6568.625100000001
$ws.Range("L138").Value = 11806.4064
$ws.Range("M138").Value = -1428.625100000001
$ws.Range("N138").Value = -22086.4064

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 295
$ws.Range("I5").Value = 295
$ws.Range("K5").Value = 295
$ws.Range("M5").Value = -183
$ws.Range("H32").Value = 4823.9614
$ws.Range("I32").Value = 5121.7915
$ws.Range("J32").Value = 1250
$ws.Range("K32").Value = 5121.7915
$ws.Range("L32").Value = 1250
$ws.Range("M32").Value = -4834.7915
$ws.Range("N32").Value = -1824
$ws.Range("H45").Value = 1384.4
$ws.Range("I45").Value = 1384.4
$ws.Range("K45").Value = 1384.4
$ws.Range("M45").Value = -1007.4
$ws.Range("H97").Value = 284.14285
$ws.Range("I97").Value = 267.6154
$ws.Range("K97").Value = 267.6154
$ws.Range("M97").Value = 228.3846
$ws.Range("H102").Value = 4637.5557
$ws.Range("I102").Value = 1956.3334
$ws.Range("K102").Value = 1956.3334
$ws.Range("M102").Value = -334.3334
$ws.Range("H132").Value = 2591.7407
$ws.Range("I132").Value = 2599.1155
$ws.Range("K132").Value = 7797.3465
$ws.Range("M132").Value = -5267.3465

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 295
$ws.Range("I4").Value = 295
$ws.Range("K4").Value = 295
$ws.Range("M4").Value = -180
$ws.Range("H86").Value = 3044.4285
$ws.Range("I86").Value = 3050.8333
$ws.Range("K86").Value = 3050.8333
$ws.Range("M86").Value = -1927.8333
$ws.Range("H89").Value = 3044.4285
$ws.Range("I89").Value = 3050.8333
$ws.Range("K89").Value = 15254.1665
$ws.Range("M89").Value = -9638.166499999999
$ws.Range("H105").Value = 1942.4117
$ws.Range("I105").Value = 1094.1818
$ws.Range("J105").Value = 3497.5
$ws.Range("K105").Value = 1094.1818
$ws.Range("L105").Value = 3497.5
$ws.Range("M105").Value = 652.8181999999999
$ws.Range("N105").Value = -6991.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 167.91667
$ws.Range("I7").Value = 142.72728
$ws.Range("K7").Value = 142.72728
$ws.Range("M7").Value = -29.72728000000001
$ws.Range("H31").Value = 2861.1765
$ws.Range("I31").Value = 1557.0667
$ws.Range("K31").Value = 1557.0667
$ws.Range("M31").Value = -1262.0667
$ws.Range("H34").Value = 2861.1765
$ws.Range("I34").Value = 1557.0667
$ws.Range("K34").Value = 1557.0667
$ws.Range("M34").Value = -1355.0667
$ws.Range("H86").Value = 8150.4443
$ws.Range("I86").Value = 7793.095
$ws.Range("K86").Value = 7793.095
$ws.Range("M86").Value = -6670.095
$ws.Range("H89").Value = 8150.4443
$ws.Range("I89").Value = 7793.095
$ws.Range("K89").Value = 38965.475
$ws.Range("M89").Value = -33349.475
$ws.Range("H132").Value = 2273.5881
$ws.Range("I132").Value = 2328.1875
$ws.Range("K132").Value = 6984.5625
$ws.Range("M132").Value = -4454.5625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 4781.8184
$ws.Range("I7").Value = 800
$ws.Range("J7").Value = 5180
$ws.Range("K7").Value = 2400
$ws.Range("L7").Value = 15540
$ws.Range("M7").Value = -2288
$ws.Range("N7").Value = -15764
$ws.Range("H10").Value = 202.5
$ws.Range("I10").Value = 69.583336
$ws.Range("K10").Value = 208.750008
$ws.Range("M10").Value = -69.75000800000001
$ws.Range("H70").Value = 3006
$ws.Range("I70").Value = 3006
$ws.Range("K70").Value = 9018
$ws.Range("M70").Value = -8703
$ws.Range("H73").Value = 3006
$ws.Range("I73").Value = 3006
$ws.Range("K73").Value = 9018
$ws.Range("M73").Value = -7926
$ws.Range("H108").Value = 12604.454
$ws.Range("I108").Value = 405.44446
$ws.Range("K108").Value = 1216.33338
$ws.Range("M108").Value = 1663.66662
$ws.Range("H132").Value = 2359.72
$ws.Range("J132").Value = 2569.7
$ws.Range("L132").Value = 23127.3
$ws.Range("N132").Value = -28187.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1425
$ws.Range("I102").Value = 1450
$ws.Range("J102").Value = 1400
$ws.Range("K102").Value = 1450
$ws.Range("L102").Value = 1400
$ws.Range("M102").Value = 172
$ws.Range("N102").Value = -4644
$ws.Range("H122").Value = 4161.353
$ws.Range("I122").Value = 3917.5
$ws.Range("K122").Value = 11752.5
$ws.Range("M122").Value = -9302.5
$ws.Range("H132").Value = 3271.7273
$ws.Range("I132").Value = 2948.95
$ws.Range("K132").Value = 8846.849999999999
$ws.Range("M132").Value = -6316.849999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3052.0293
$ws.Range("I22").Value = 2955.8572
$ws.Range("J22").Value = 3119.35
$ws.Range("K22").Value = 2955.8572
$ws.Range("L22").Value = 3119.35
$ws.Range("M22").Value = -2660.8572
$ws.Range("N22").Value = -3709.35
$ws.Range("H27").Value = 3052.0293
$ws.Range("I27").Value = 2955.8572
$ws.Range("J27").Value = 3119.35
$ws.Range("K27").Value = 2955.8572
$ws.Range("L27").Value = 3119.35
$ws.Range("M27").Value = -2848.8572
$ws.Range("N27").Value = -3333.35
$ws.Range("H46").Value = 2454.862
$ws.Range("I46").Value = 1230.4546
$ws.Range("K46").Value = 1230.4546
$ws.Range("M46").Value = -1042.4546
$ws.Range("H100").Value = 6446.909
$ws.Range("I100").Value = 4000
$ws.Range("K100").Value = 4000
$ws.Range("M100").Value = -3459
$ws.Range("H136").Value = 3640.5
$ws.Range("I136").Value = 4159.8335
$ws.Range("J136").Value = 3195.3572
$ws.Range("K136").Value = 12479.5005
$ws.Range("L136").Value = 9586.071599999999
$ws.Range("M136").Value = -9929.500499999998
$ws.Range("N136").Value = -14686.0716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 28489.75
$ws.Range("J74").Value = 24319.666
$ws.Range("L74").Value = 24319.666
$ws.Range("N74").Value = -26191.666
$ws.Range("H77").Value = 28489.75
$ws.Range("J77").Value = 24319.666
$ws.Range("L77").Value = 72958.99800000001
$ws.Range("N77").Value = -82318.99800000001
